$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new rule row above the current row 24 ("Assign Case Number") ---
# Insert blank cells (shifting existing rows down) in A24:F24, then copy the
# formatting that now lives at row 25 (the row that used to be row 24) back
# up onto the freshly inserted row 24, so the new row matches its neighbours.
$ws.Range("A24:F24").Insert()
$ws.Range("A25:F25").Copy()
$ws.Range("A24:F24").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match row height of the row that the new row was cloned from.
$ws.Rows.Item(24).RowHeight = $ws.Rows.Item(25).RowHeight

# Populate the new "Set Response Due Date" rule.
$ws.Range("B24").Value = "Set Response Due Date"
$ws.Range("C24").Value = "responseDueDate == null"
$ws.Range("D24").Value = "setResponseDueDate, java.time.LocalDate.now().plusDays(10)"

# --- Remove the obsolete "setDueDate, null" action on the (now) last rule row ---
# That rule used to live on row 33 ("Set Due Date Release Queue"); after the
# insertion above it has shifted down to row 34.
$ws.Range("D34").ClearContents()
